# Swedish steel and simplex method
# --------------------------------
# The sheet models an LP (transportation-style) problem. This edit reworks
# sheet1 ("Foglio1") from a plain A:H constraint table into a simplex-method
# tableau: the redundant "all -1" row is dropped, ten slack/identity columns
# (H:Q) are inserted before the RHS column, and the RHS values move from
# column H out to column R.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# 1. Drop the old row 3 (the redundant "-1,-1,...,-1000" row). Everything
#    below shifts up one row automatically (A1:H13 -> A1:H12), formulas in
#    A8/B8 included.
$ws.Rows("3:3").Delete()

# 2. Move the RHS column (currently H, holding 1000,-6.5,7.5,-30,... ,250)
#    out to its final home in column R, to make room for the identity
#    (slack-variable) block.
$ws.Range("H1:H12").Cut($ws.Range("R1:R12"))

# The cell that used to be H5 (value 7.5) carried a quote-prefix cell style
# (s="1"); after the shift/cut that style stayed behind on H4 (now empty)
# instead of following the value to R4. Move the formatting across, then
# reset the now-empty H4 back to the plain/default format.
$ws.Range("H4").Copy()
$ws.Range("R4").PasteSpecial(-4122)
$ws.Range("A4").Copy()
$ws.Range("H4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 3. Fill H2:Q12 with the simplex identity block: row 2 (the objective /
#    first data row) is all zeros, and rows 3-12 each carry a single 1 that
#    walks across H..Q (a 10x10 identity matrix for the slack variables).
$identity = New-Object 'object[,]' 11,10
for ($i = 0; $i -lt 11; $i++) {
    for ($j = 0; $j -lt 10; $j++) {
        $identity[$i, $j] = 0
    }
}
for ($i = 1; $i -lt 11; $i++) {
    $identity[$i, $i - 1] = 1
}
$ws.Range("H2:Q12").Value = $identity

# 4. Restore the view: the saved selection/scroll now point at the new
#    bottom-right corner of the tableau.
$ws.Range("P13").Select()
$excel.ActiveWindow.ScrollColumn = 2
